# Updated cryptos list (price + volume refresh, plus two row reorderings)
# Leading "'" forces text so prices like "37.501.62" / "2.60" are not
# auto-converted into numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''37.501.62'
$ws.Range('E2').Value = '  +2.54%  '

$ws.Range('D3').Value = '''2.077.97'
$ws.Range('E3').Value = '  +3.66%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').Value = '''235.07'
$ws.Range('E5').Value = '  -0.84%  '

$ws.Range('D6').Value = '''0.622'
$ws.Range('E6').Value = '  +3.76%  '

$ws.Range('E7').Value = '  +6.02%  '

$ws.Range('E9').Value = '  +3.86%  '

$ws.Range('D10').Value = '''59.05'
$ws.Range('E10').Value = '  +1.29%  '

$ws.Range('E11').Value = '  +2.04%  '

$ws.Range('E12').Value = '  +3.82%  '

$ws.Range('D13').Value = '''2.385.01'
$ws.Range('E13').Value = '  +3.74%  '

$ws.Range('D14').Value = '''14.58'
$ws.Range('E14').Value = '  +2.72%  '

$ws.Range('D15').Value = '''21.15'
$ws.Range('E15').Value = '  +4.50%  '

$ws.Range('E16').Value = '  +3.19%  '

$ws.Range('D17').Value = '''5.19'
$ws.Range('E17').Value = '  +1.79%  '

$ws.Range('D18').Value = '''2.075.14'
$ws.Range('E18').Value = '  +3.49%  '

$ws.Range('D19').Value = '''37.448.45'
$ws.Range('E19').Value = '  +2.62%  '

$ws.Range('D20').Value = '''6.26'
$ws.Range('E20').Value = '  +18.11%  '

$ws.Range('D21').Value = '''69.92'
$ws.Range('E21').Value = '  +3.10%  '

$ws.Range('E22').Value = '  +1.68%  '

$ws.Range('D23').Value = '''226.60'
$ws.Range('E23').Value = '  +2.19%  '

$ws.Range('E24').Value = '  -0.21%  '

$ws.Range('D25').Value = '''2.47'
$ws.Range('E25').Value = '  +2.78%  '

$ws.Range('D26').Value = '''2.38'
$ws.Range('E26').Value = '  +0.46%  '

$ws.Range('D27').Value = '''166.97'
$ws.Range('E27').Value = '  +3.21%  '

$ws.Range('D28').Value = '''1.51'
$ws.Range('E28').Value = '  +11.19%  '

$ws.Range('E29').Value = '  +4.36%  '

$ws.Range('D30').Value = '''19.27'
$ws.Range('E30').Value = '  +2.39%  '

$ws.Range('E31').Value = '  -0.78%  '

$ws.Range('E32').Value = '  +0.93%  '

$ws.Range('E33').Value = '  +3.75%  '

$ws.Range('B34').Value = '''LidoDAOToken'
$ws.Range('C34').Value = '''https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = '''2.60'
$ws.Range('E34').Value = '  +6.52%  '

$ws.Range('B35').Value = '''Hedera'
$ws.Range('C35').Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '''0.0625'
$ws.Range('E35').Value = '  +3.45%  '

$ws.Range('E36').Value = '  +7.55%  '

$ws.Range('E37').Value = '  -0.08%  '

$ws.Range('D38').Value = '''3.35'
$ws.Range('E38').Value = '  +0.25%  '

$ws.Range('D39').Value = '''5.88'
$ws.Range('E39').Value = '  +3.04%  '

$ws.Range('E40').Value = '  +0.40%  '

$ws.Range('D41').Value = '''4.61'
$ws.Range('E41').Value = '  +21.90%  '

$ws.Range('E42').Value = '  -1.07%  '

$ws.Range('D43').Value = '''0.0957'
$ws.Range('E43').Value = '  +3.50%  '

$ws.Range('D44').Value = '''1.476.41'
$ws.Range('E44').Value = '  +1.70%  '

$ws.Range('B45').Value = '''TrustWalletToken'
$ws.Range('C45').Value = '''https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').Value = '''1.18'
$ws.Range('E45').Value = '  +6.85%  '

$ws.Range('B46').Value = '''Aave'
$ws.Range('C46').Value = '''https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '''96.03'
$ws.Range('E46').Value = '  +6.28%  '

$ws.Range('E47').Value = '  +4.73%  '

$ws.Range('D48').Value = '''15.86'
$ws.Range('E48').Value = '  +3.96%  '

$ws.Range('E49').Value = '  +3.71%  '

$ws.Range('D50').Value = '''7.30'
$ws.Range('E50').Value = '  +6.58%  '

$ws.Range('E51').Value = '  +1.74%  '

